# Adding new data sets
# The workbook already contains all the rows of data; this edit applies a
# filtered view over it: filter the existing AutoFilter down to
# Genotype=StAug2, Inoculate=WSM, SoilConc in {high, low}, Root/Shoot=Shoot,
# which hides the rows that don't match, grows the AutoFilter range to
# include the last data row (L79), and leaves the selection/pane parked on
# the range the analyst was last looking at (L21:L41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate so the AutoFilter is rebuilt (and its range
# re-measured against the full used range) rather than just resized in
# place.
$ws.AutoFilterMode = $false

$used = $ws.UsedRange
$filterRange = $ws.Range("A1:L79")

# Genotype (column B, colId 1 in the filter) = StAug2
$filterRange.AutoFilter(2, @("StAug2"), 7)
# Inoculate (column G, colId 6) = WSM
$filterRange.AutoFilter(7, @("WSM"), 7)
# SoilConc (column H, colId 7) = high or low
$filterRange.AutoFilter(8, @("high", "low"), 7)
# Root/Shoot (column K, colId 10) = Shoot
$filterRange.AutoFilter(11, @("Shoot"), 7)

# The hidden "_FilterDatabase" defined name Excel keeps in sync with the
# AutoFilter range doesn't auto-resize here, so push it out explicitly.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$L`$79"

# Park the selection/pane where the analyst left it.
$ws.Range("L21:L41").Select()
